$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("shulConfigeration")

$ws.Range("D89").NumberFormat = "@"
$ws.Range("D89").Value = "905"
$ws.Range("D90").NumberFormat = "@"
$ws.Range("D90").Value = "915"
$ws.Range("D91").NumberFormat = "@"
$ws.Range("D91").Value = "925"

$ws.Range("G113").Select()
